$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 70: add JENIS KERTAS (C) and UKURAN (D), update price (E)
$ws.Range("C70").Value = "HVS (GOLD)"
$ws.Range("D70").Value = "A7 (7,4X10,5CM)"
$ws.Range("E70").Value = 8500

# Row 71: new row
$ws.Range("B71").Value = "PAKET MURAH ALQURAN AL AQEEL MUSHAF NON TERJEMAHAN | SURABAYA | al quran Wakaf/Shodaqoh hadiah hampers islami"
$ws.Range("C71").Value = "HVS (GOLD)"
$ws.Range("D71").Value = "A5 (14,5X20CM)"
$ws.Range("E71").Value = 21000

# Row 72: new row
$ws.Range("B72").Value = "PAKET MURAH ALQURAN AL AQEEL MUSHAF NON TERJEMAHAN | SURABAYA | al quran Wakaf/Shodaqoh hadiah hampers islami"
$ws.Range("C72").Value = "HVS"
$ws.Range("D72").Value = "A6 (10,5X14,5CM)"
$ws.Range("E72").Value = 12500

# Row 73: new row
$ws.Range("B73").Value = "PAKET MURAH ALQURAN AL AQEEL MUSHAF NON TERJEMAHAN | SURABAYA | al quran Wakaf/Shodaqoh hadiah hampers islami"
$ws.Range("C73").Value = "KORAN"
$ws.Range("D73").Value = "A5 (14,5X20CM)"
$ws.Range("E73").Value = 14500
